## Add a new "Skills_Matrix" worksheet containing a Brew/Filter/Pack/Ship/Manage
## skills matrix for five people, and turn that range into a table ("Skills_Mtx").

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing sheet (becomes the active tab).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Skills_Matrix"

# Names down column A first (row by row), then the skill headers across row 1,
# then the Column1 header for the name column.
$ws.Range("A2").Value = "Alfred"
$ws.Range("A3").Value = "Bill"
$ws.Range("A4").Value = "Chris"
$ws.Range("A5").Value = "Dante"
$ws.Range("A6").Value = "Edgar"

$ws.Range("B1").Value = "Brew"
$ws.Range("C1").Value = "Filter"
$ws.Range("D1").Value = "Pack"
$ws.Range("E1").Value = "Ship"
$ws.Range("F1").Value = "Manage"

$ws.Range("A1").Value = "Column1"

# 0/1 matrix of who covers which skill.
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0

# Turn the range into a table, matching the header names already in row 1.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:F6"), $null, 1)
$tbl.Name = "Skills_Mtx"
$tbl.TableStyle = "TableStyleLight1"

# Size columns to fit their contents and land on the cell the author left selected.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Range("D3").Select() | Out-Null
